$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(20).Insert()

$ws.Cells.Item(20, 1).Value = 7
$ws.Cells.Item(20, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(20, 3).Value = "Ñuble"
$ws.Cells.Item(20, 4).Value = 44845
$ws.Cells.Item(20, 5).Value = 16
$ws.Cells.Item(20, 6).Value = 300000000
$ws.Cells.Item(20, 7).Value = "Espárragos"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 400
$ws.Cells.Item(20, 11).Value = 1300
$ws.Cells.Item(20, 12).Value = 1500
$ws.Cells.Item(20, 13).Value = 1400
$ws.Cells.Item(20, 14).Value = "$/kilo"
$ws.Cells.Item(20, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(20, 16).Value = 1400
$ws.Cells.Item(20, 17).Value = 1
$ws.Cells.Item(20, 18).Value = "Hortaliza"

Write-Host "Done"
